$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.07901576910072249
$ws.Range("J2").Value = 0.07901576910072249
$ws.Range("M2").Value = 0.448453
$ws.Range("N2").Value = 1.345359
$ws.Range("O2").Value = 0.3700410539697698
$ws.Range("P2").Value = 0.3700410539697698
$ws.Range("Q2").Value = 0.04412045046766667
$ws.Range("R2").Value = 0.3970840542090001
$ws.Range("S2").Value = 0.02923907847826332
$ws.Range("T2").Value = 0.02923907847826332
# Row 3
$ws.Range("I3").Value = 0.07901576910072249
$ws.Range("J3").Value = 0.07901576910072249
$ws.Range("O3").Value = 0.2028092511432455
$ws.Range("P3").Value = 0.2028092511432455
$ws.Range("S3").Value = 0.01602512895982513
$ws.Range("T3").Value = 0.01602512895982513
# Row 4
$ws.Range("I4").Value = 0.07901576910072249
$ws.Range("J4").Value = 0.07901576910072249
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03878766666666666
$ws.Range("N4").Value = 0.116363
$ws.Range("O4").Value = 0.03200564842773142
$ws.Range("P4").Value = 0.03200564842773142
$ws.Range("Q4").Value = 0.003816072868111111
$ws.Range("R4").Value = 0.034344655813
$ws.Range("S4").Value = 0.002528950926084528
$ws.Range("T4").Value = 0.002528950926084528
# Row 5
$ws.Range("I5").Value = 0.07901576910072249
$ws.Range("J5").Value = 0.07901576910072249
$ws.Range("M5").Value = 0.188246
$ws.Range("N5").Value = 0.564738
$ws.Range("O5").Value = 0.1553312125141169
$ws.Range("P5").Value = 0.1553312125141169
$ws.Range("Q5").Value = 0.01852033171533333
$ws.Range("R5").Value = 0.166682985438
$ws.Range("S5").Value = 0.01227361522215072
$ws.Range("T5").Value = 0.01227361522215072
# Row 6
$ws.Range("I6").Value = 0.07901576910072249
$ws.Range("J6").Value = 0.07901576910072249
$ws.Range("M6").Value = 0.1785163333333334
$ws.Range("N6").Value = 0.5355490000000001
$ws.Range("O6").Value = 0.1473027767402279
$ws.Range("P6").Value = 0.1473027767402279
$ws.Range("Q6").Value = 0.01756309143322222
$ws.Range("R6").Value = 0.158067822899
$ws.Range("S6").Value = 0.01163924219480112
$ws.Range("T6").Value = 0.01163924219480112
# Row 7
$ws.Range("I7").Value = 0.07901576910072249
$ws.Range("J7").Value = 0.07901576910072249
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.112113
$ws.Range("N7").Value = 0.336339
$ws.Range("O7").Value = 0.09251005720490843
$ws.Range("P7").Value = 0.09251005720490843
$ws.Range("Q7").Value = 0.011030088021
$ws.Range("R7").Value = 0.09927079218899999
$ws.Range("S7").Value = 0.007309753319597673
$ws.Range("T7").Value = 0.007309753319597673
# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1928733333333333
$ws.Range("H8").Value = 0.5786199999999999
$ws.Range("I8").Value = 0.1549041145619023
$ws.Range("J8").Value = 0.1549041145619023
$ws.Range("M8").Value = 0.448453
$ws.Range("N8").Value = 1.345359
$ws.Range("O8").Value = 0.3700410539697698
$ws.Range("P8").Value = 0.3700410539697698
$ws.Range("Q8").Value = 0.08649462495333334
$ws.Range("R8").Value = 0.7784516245799999
$ws.Range("S8").Value = 0.05732088181674031
$ws.Range("T8").Value = 0.05732088181674031
# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1928733333333333
$ws.Range("H9").Value = 0.5786199999999999
$ws.Range("I9").Value = 0.1549041145619023
$ws.Range("J9").Value = 0.1549041145619023
$ws.Range("O9").Value = 0.2028092511432455
$ws.Range("P9").Value = 0.2028092511432455
$ws.Range("Q9").Value = 0.04740530794222222
$ws.Range("R9").Value = 0.4266477714799999
$ws.Range("S9").Value = 0.03141598747330692
$ws.Range("T9").Value = 0.03141598747330692
# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1928733333333333
$ws.Range("H10").Value = 0.5786199999999999
$ws.Range("I10").Value = 0.1549041145619023
$ws.Range("J10").Value = 0.1549041145619023
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03878766666666666
$ws.Range("N10").Value = 0.116363
$ws.Range("O10").Value = 0.03200564842773142
$ws.Range("P10").Value = 0.03200564842773142
$ws.Range("Q10").Value = 0.007481106562222221
$ws.Range("R10").Value = 0.06732995905999999
$ws.Range("S10").Value = 0.004957806630677278
$ws.Range("T10").Value = 0.004957806630677278
# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1928733333333333
$ws.Range("H11").Value = 0.5786199999999999
$ws.Range("I11").Value = 0.1549041145619023
$ws.Range("J11").Value = 0.1549041145619023
$ws.Range("M11").Value = 0.188246
$ws.Range("N11").Value = 0.564738
$ws.Range("O11").Value = 0.1553312125141169
$ws.Range("P11").Value = 0.1553312125141169
$ws.Range("Q11").Value = 0.03630763350666666
$ws.Range("R11").Value = 0.3267687015599999
$ws.Range("S11").Value = 0.02406144393832597
$ws.Range("T11").Value = 0.02406144393832596
# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1928733333333333
$ws.Range("H12").Value = 0.5786199999999999
$ws.Range("I12").Value = 0.1549041145619023
$ws.Range("J12").Value = 0.1549041145619023
$ws.Range("M12").Value = 0.1785163333333334
$ws.Range("N12").Value = 0.5355490000000001
$ws.Range("O12").Value = 0.1473027767402279
$ws.Range("P12").Value = 0.1473027767402279
$ws.Range("Q12").Value = 0.03443104026444445
$ws.Range("R12").Value = 0.30987936238
$ws.Range("S12").Value = 0.02281780620345459
$ws.Range("T12").Value = 0.02281780620345459
# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1928733333333333
$ws.Range("H13").Value = 0.5786199999999999
$ws.Range("I13").Value = 0.1549041145619023
$ws.Range("J13").Value = 0.1549041145619023
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.112113
$ws.Range("N13").Value = 0.336339
$ws.Range("O13").Value = 0.09251005720490843
$ws.Range("P13").Value = 0.09251005720490843
$ws.Range("Q13").Value = 0.02162360802
$ws.Range("R13").Value = 0.19461247218
$ws.Range("S13").Value = 0.01433018849939727
$ws.Range("T13").Value = 0.01433018849939727
# Row 14
$ws.Range("G14").Value = 0.9538573333333332
$ws.Range("H14").Value = 2.861572
$ws.Range("I14").Value = 0.7660801163373752
$ws.Range("J14").Value = 0.7660801163373753
$ws.Range("M14").Value = 0.448453
$ws.Range("N14").Value = 1.345359
$ws.Range("O14").Value = 0.3700410539697698
$ws.Range("P14").Value = 0.3700410539697698
$ws.Range("Q14").Value = 0.4277601827053333
$ws.Range("R14").Value = 3.849841644348
$ws.Range("S14").Value = 0.2834810936747662
$ws.Range("T14").Value = 0.2834810936747662
# Row 15
$ws.Range("G15").Value = 0.9538573333333332
$ws.Range("H15").Value = 2.861572
$ws.Range("I15").Value = 0.7660801163373752
$ws.Range("J15").Value = 0.7660801163373753
$ws.Range("O15").Value = 0.2028092511432455
$ws.Range("P15").Value = 0.2028092511432455
$ws.Range("Q15").Value = 0.2344435067208888
$ws.Range("R15").Value = 2.109991560488
$ws.Range("S15").Value = 0.1553681347101134
$ws.Range("T15").Value = 0.1553681347101135
# Row 16
$ws.Range("G16").Value = 0.9538573333333332
$ws.Range("H16").Value = 2.861572
$ws.Range("I16").Value = 0.7660801163373752
$ws.Range("J16").Value = 0.7660801163373753
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03878766666666666
$ws.Range("N16").Value = 0.116363
$ws.Range("O16").Value = 0.03200564842773142
$ws.Range("P16").Value = 0.03200564842773142
$ws.Range("Q16").Value = 0.03699790029288888
$ws.Range("R16").Value = 0.3329811026359999
$ws.Range("S16").Value = 0.02451889087096962
$ws.Range("T16").Value = 0.02451889087096962
# Row 17
$ws.Range("G17").Value = 0.9538573333333332
$ws.Range("H17").Value = 2.861572
$ws.Range("I17").Value = 0.7660801163373752
$ws.Range("J17").Value = 0.7660801163373753
$ws.Range("M17").Value = 0.188246
$ws.Range("N17").Value = 0.564738
$ws.Range("O17").Value = 0.1553312125141169
$ws.Range("P17").Value = 0.1553312125141169
$ws.Range("Q17").Value = 0.1795598275706666
$ws.Range("R17").Value = 1.616038448136
$ws.Range("S17").Value = 0.1189961533536402
$ws.Range("T17").Value = 0.1189961533536402
# Row 18
$ws.Range("G18").Value = 0.9538573333333332
$ws.Range("H18").Value = 2.861572
$ws.Range("I18").Value = 0.7660801163373752
$ws.Range("J18").Value = 0.7660801163373753
$ws.Range("M18").Value = 0.1785163333333334
$ws.Range("N18").Value = 0.5355490000000001
$ws.Range("O18").Value = 0.1473027767402279
$ws.Range("P18").Value = 0.1473027767402279
$ws.Range("Q18").Value = 0.1702791136697778
$ws.Range("R18").Value = 1.532512023028
$ws.Range("S18").Value = 0.1128457283419722
$ws.Range("T18").Value = 0.1128457283419722
# Row 19
$ws.Range("G19").Value = 0.9538573333333332
$ws.Range("H19").Value = 2.861572
$ws.Range("I19").Value = 0.7660801163373752
$ws.Range("J19").Value = 0.7660801163373753
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.112113
$ws.Range("N19").Value = 0.336339
$ws.Range("O19").Value = 0.09251005720490843
$ws.Range("P19").Value = 0.09251005720490843
$ws.Range("Q19").Value = 0.106939807212
$ws.Range("R19").Value = 0.9624582649079999
$ws.Range("S19").Value = 0.07087011538591348
$ws.Range("T19").Value = 0.07087011538591349

Write-Host "Applied 220 cell updates"
